$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("福龙马", "三六零", "平潭发展")
    3  = @("平潭发展", "福龙马", "山子高科")
    4  = @("三六零", "万向钱潮", "福龙马")
    5  = @("大众公用", "平潭发展", "三六零")
    6  = @("山子高科", "闻泰科技", "万向钱潮")
    7  = @("万向钱潮", "多氟多", "多氟多")
    8  = @("多氟多", "天际股份", "东方精工")
    9  = @("天际股份", "大众公用", "天际股份")
    10 = @("闻泰科技", "粤 传 媒", "大众公用")
    11 = @("包钢股份", "山子高科", "时空科技")
    12 = @("东方精工", "包钢股份", "粤传媒")
    13 = @("粤 传 媒", "东方精工", "三花智控")
    14 = @("东方明珠", "上海电气", "闻泰科技")
    15 = @("中国核建", "中国核建", "神州信息")
    16 = @("上海电气", "永鼎股份", "中国核建")
    17 = @("贵州茅台", "海马汽车", "和而泰")
    18 = @("时空科技", "亚太药业", "工业富联")
    19 = @("海峡创新", "鼎胜新材", "上海电气")
    20 = @("海马汽车", "贵州茅台", "海峡创新")
    21 = @("亚太药业", "东方财富", "利欧股份")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
}
